$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "yt5"
$ws.Range("B1").Value = 555
$ws.Range("C1").Value = 555
$ws.Range("E1").Value = "5ttt"
$ws.Range("B4").Value = "gg"
$ws.Range("C3").Value = "g"

$ws.Range("C3").Select()
